$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the numeric-looking stat columns (C:G) as text, matching the
# original workbook where every data cell is stored as a shared string.
$ws.Range("C2:G19").NumberFormat = "@"

$ws.Range("B2").Value = "PSV"
$ws.Range("C2").Value = "1.0"
$ws.Range("D2").Value = "7.5"
$ws.Range("E2").Value = "100%"
$ws.Range("F2").Value = "75%"
$ws.Range("G2").Value = "3.83"
$ws.Range("B3").Value = "Feyenoord"
$ws.Range("C3").Value = "0.9"
$ws.Range("D3").Value = "7.4"
$ws.Range("E3").Value = "79%"
$ws.Range("F3").Value = "67%"
$ws.Range("G3").Value = "3.29"
$ws.Range("B4").Value = "Twente"
$ws.Range("C4").Value = "1.6"
$ws.Range("D4").Value = "6.5"
$ws.Range("E4").Value = "75%"
$ws.Range("F4").Value = "67%"
$ws.Range("G4").Value = "2.92"
$ws.Range("B5").Value = "AZ"
$ws.Range("C5").Value = "1.3"
$ws.Range("D5").Value = "6.3"
$ws.Range("E5").Value = "88%"
$ws.Range("F5").Value = "58%"
$ws.Range("G5").Value = "2.83"
$ws.Range("B6").Value = "Ajax"
$ws.Range("C6").Value = "1.2"
$ws.Range("D6").Value = "6.1"
$ws.Range("E6").Value = "96%"
$ws.Range("F6").Value = "79%"
$ws.Range("G6").Value = "4.00"
$ws.Range("B7").Value = "G. A. Eagle"
$ws.Range("C7").Value = "1.4"
$ws.Range("D7").Value = "5.1"
$ws.Range("E7").Value = "83%"
$ws.Range("F7").Value = "50%"
$ws.Range("G7").Value = "2.92"
$ws.Range("B8").Value = "Nijmegen"
$ws.Range("C8").Value = "1.7"
$ws.Range("D8").Value = "4.3"
$ws.Range("E8").Value = "96%"
$ws.Range("F8").Value = "71%"
$ws.Range("G8").Value = "3.63"
$ws.Range("B9").Value = "Utrecht"
$ws.Range("C9").Value = "1.5"
$ws.Range("D9").Value = "5.8"
$ws.Range("E9").Value = "75%"
$ws.Range("F9").Value = "33%"
$ws.Range("G9").Value = "2.54"
$ws.Range("B10").Value = "Heerenveen"
$ws.Range("C10").Value = "1.3"
$ws.Range("D10").Value = "5.6"
$ws.Range("E10").Value = "96%"
$ws.Range("F10").Value = "75%"
$ws.Range("G10").Value = "3.50"
$ws.Range("B11").Value = "Sparta Rotterdam"
$ws.Range("C11").Value = "1.5"
$ws.Range("D11").Value = "5.4"
$ws.Range("E11").Value = "88%"
$ws.Range("F11").Value = "46%"
$ws.Range("G11").Value = "2.67"
$ws.Range("B12").Value = "Fortuna Sittard"
$ws.Range("C12").Value = "1.8"
$ws.Range("D12").Value = "5.1"
$ws.Range("E12").Value = "71%"
$ws.Range("F12").Value = "63%"
$ws.Range("G12").Value = "2.88"
$ws.Range("B13").Value = "Almere City"
$ws.Range("C13").Value = "2.2"
$ws.Range("D13").Value = "4.5"
$ws.Range("E13").Value = "71%"
$ws.Range("F13").Value = "54%"
$ws.Range("G13").Value = "2.71"
$ws.Range("B14").Value = "Zwolle"
$ws.Range("C14").Value = "2.0"
$ws.Range("D14").Value = "3.5"
$ws.Range("E14").Value = "88%"
$ws.Range("F14").Value = "58%"
$ws.Range("G14").Value = "3.13"
$ws.Range("B15").Value = "Heracles Almelo"
$ws.Range("C15").Value = "2.0"
$ws.Range("D15").Value = "4.5"
$ws.Range("E15").Value = "92%"
$ws.Range("F15").Value = "75%"
$ws.Range("G15").Value = "3.58"
$ws.Range("B16").Value = "Excelsior"
$ws.Range("C16").Value = "1.5"
$ws.Range("D16").Value = "3.2"
$ws.Range("E16").Value = "92%"
$ws.Range("F16").Value = "79%"
$ws.Range("G16").Value = "3.75"
$ws.Range("B17").Value = "Waalwijk"
$ws.Range("C17").Value = "1.6"
$ws.Range("D17").Value = "3.8"
$ws.Range("E17").Value = "71%"
$ws.Range("F17").Value = "50%"
$ws.Range("G17").Value = "2.54"
$ws.Range("B18").Value = "Vitesse"
$ws.Range("C18").Value = "1.7"
$ws.Range("D18").Value = "5.2"
$ws.Range("E18").Value = "79%"
$ws.Range("F18").Value = "54%"
$ws.Range("G18").Value = "2.83"
$ws.Range("B19").Value = "FC Volendam"
$ws.Range("C19").Value = "1.5"
$ws.Range("D19").Value = "4.3"
$ws.Range("E19").Value = "88%"
$ws.Range("F19").Value = "71%"
$ws.Range("G19").Value = "3.63"
